$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.736660666666666
$ws.Range("H2").Value = 23.209982
$ws.Range("M2").Value = 0.7861523333333333
$ws.Range("N2").Value = 2.358457
$ws.Range("O2").Value = 0.01668797875153133
$ws.Range("P2").Value = 0.01668797875153133
$ws.Range("Q2").Value = 6.082193835308222
$ws.Range("R2").Value = 54.739744517774
$ws.Range("S2").Value = 0.01668797875153133
$ws.Range("T2").Value = 0.01668797875153133

$ws.Range("G3").Value = 7.736660666666666
$ws.Range("H3").Value = 23.209982
$ws.Range("O3").Value = 0.5736784050900728
$ws.Range("P3").Value = 0.5736784050900727
$ws.Range("Q3").Value = 209.086031977846
$ws.Range("R3").Value = 1881.774287800614
$ws.Range("S3").Value = 0.5736784050900728
$ws.Range("T3").Value = 0.5736784050900727

$ws.Range("G4").Value = 7.736660666666666
$ws.Range("H4").Value = 23.209982
$ws.Range("M4").Value = 19.22475933333333
$ws.Range("N4").Value = 57.674278
$ws.Range("O4").Value = 0.4080918692916219
$ws.Range("P4").Value = 0.4080918692916219
$ws.Range("Q4").Value = 148.7354393603329
$ws.Range("R4").Value = 1338.618954242996
$ws.Range("S4").Value = 0.4080918692916219
$ws.Range("T4").Value = 0.4080918692916219

$ws.Range("G5").Value = 7.736660666666666
$ws.Range("H5").Value = 23.209982
$ws.Range("M5").Value = 0.07263
$ws.Range("N5").Value = 0.21789
$ws.Range("O5").Value = 0.00154174686677398
$ws.Range("P5").Value = 0.00154174686677398
$ws.Range("Q5").Value = 0.56191366422
$ws.Range("R5").Value = 5.05722297798
$ws.Range("S5").Value = 0.00154174686677398
$ws.Range("T5").Value = 0.00154174686677398
